$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.061.01'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '3.574.28'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '602.39'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '135.88'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.52%  '
$ws.Range('D7').Value = '3.573.70'
$ws.Range('E7').Value = '  +2.71%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +1.23%  '
$ws.Range('E10').Value = '  +0.72%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.95'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.70%  '
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('D13').Value = '4.180.54'
$ws.Range('E13').Value = '  +2.67%  '
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').Value = '3.573.92'
$ws.Range('E15').Value = '  +3.48%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '27.19'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.19%  '
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').Value = '65.148.60'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('E19').Value = '  +4.31%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.43'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.14%  '
$ws.Range('E21').Value = '  +1.65%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '388.69'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  +4.77%  '
$ws.Range('D24').Value = '3.717.96'
$ws.Range('E24').Value = '  +2.78%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '74.26'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.53%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000117'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +7.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.74'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.62%  '
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.30'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +4.00%  '
$ws.Range('E31').Value = '  +2.30%  '
$ws.Range('E32').Value = '  +21.39%  '
$ws.Range('D33').Value = '3.582.03'
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '24.04'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +4.55%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +1.91%  '
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '169.25'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('E39').Value = '  +5.67%  '
$ws.Range('E40').Value = '  +5.92%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0809'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.36%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '27.15'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +12.85%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.827'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '42.69'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.86%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.49'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.10%  '
$ws.Range('E47').Value = '  +6.23%  '
$ws.Range('E48').Value = '  +2.64%  '
$ws.Range('D49').Value = '2.508.43'
$ws.Range('E49').Value = '  +12.73%  '
$ws.Range('E50').Value = '  +4.33%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.39'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +11.56%  '
